$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the prediction/result of the 3rd trade as a new row in the sheet.
# Numeric columns can be assigned directly.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 0.5568974771873323
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 0.8
$ws.Range("H4").Value = 0.05
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 1

# Text-like columns (dates stored as text, percentages stored as text)
# must be entered as formulas returning a string and then converted to a
# static value, otherwise Excel's value parser reinterprets them as
# numbers/percentages and changes their cell formatting.
function Set-TextValue($cell, $text) {
    $cell.Formula = "=""$text"""
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

Set-TextValue $ws.Range("C4") "20200101"
Set-TextValue $ws.Range("D4") "20210130"
Set-TextValue $ws.Range("I4") "49%"

$excel.CutCopyMode = 0
